$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2021-06-29", "overview", "K02000001", "United Kingdom", 4775301, 20479, 23, 128126),
    @("2021-06-30", "overview", "K02000001", "United Kingdom", 4800907, 26068, 14, 128140),
    @("2021-07-01", "overview", "K02000001", "United Kingdom", 4828463, 27989, 22, 128162),
    @("2021-07-02", "overview", "K02000001", "United Kingdom", 4855169, 27125, 27, 128189),
    @("2021-07-03", "overview", "K02000001", "United Kingdom", 4879616, 24885, 18, 128207),
    @("2021-07-04", "overview", "K02000001", "United Kingdom", 4903434, 24248, 15, 128222)
)

$startRow = 322
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]

    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $values[0]
    $dateCell.Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
    $ws.Cells.Item($row, 7).Value = $values[6]
    $ws.Cells.Item($row, 8).Value = $values[7]
}
